$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 156, shifting existing rows 156:285 down to 157:286
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row 156 with the new data point
$ws.Cells.Item(156, 1).Value = 10
$ws.Cells.Item(156, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(156, 3).Value = "La Araucanía"
$ws.Cells.Item(156, 4).Value = 45040
$ws.Cells.Item(156, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(156, 5).Value = 9
$ws.Cells.Item(156, 6).Value = 100112005
$ws.Cells.Item(156, 7).Value = "Puerro"
$ws.Cells.Item(156, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 50
$ws.Cells.Item(156, 11).Value = 11000
$ws.Cells.Item(156, 12).Value = 11000
$ws.Cells.Item(156, 13).Value = 11000
$ws.Cells.Item(156, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(156, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(156, 16).Value = 917
$ws.Cells.Item(156, 17).Value = 12
$ws.Cells.Item(156, 18).Value = "Hortaliza"
